# Create function to assess dumpsite scores
# Rework the "dumpsite_weight" lookup table: rename the "impact" header to
# "d_impact", flip the score_influence weights to negative values, and add
# a new "none" / 0 row so an absent impact contributes nothing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dumpsite_weight")

$ws.Range("A1").Value = "d_impact"

$ws.Range("B2").Value = -3
$ws.Range("B3").Value = -2
$ws.Range("B4").Value = -1

$ws.Range("A5").Value = "none"
$ws.Range("B5").Value = 0

# connectivity_score had been the active tab before this edit; park its
# selection on the cell it was left on, but don't activate it.
$ws4 = $wb.Worksheets.Item("connectivity_score")
$ws4.Range("H11").Select()

# Leave the user where they were working: on the dumpsite_weight sheet,
# with the selection sitting just below the new row. Activating this
# sheet last makes it the new active tab.
$ws.Activate()
$ws.Range("B6").Select()
